# Intro to Apache Spark Schedule.xlsx - schedule update
# - Removes the last two rows of the schedule table (old rows 19 & 20),
#   shifting everything up by two rows.
# - Rewrites the TIME / TOPIC columns for the (now 18-row) schedule table
#   to reflect the new timing and renamed/re-ordered exercises.
# - Row 17 (previously the highlighted "Break" row) becomes a normal,
#   un-highlighted row, so its fill/number-format is reset to match the
#   other plain rows.
# - Updates the selected cell and restores it after the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the two trailing rows of the table (19 & 20) ------------------
# This shifts rows 21.. (none exist) up and leaves rows 1-18 with whatever
# formatting rows 1-18 already had (merges / styles are preserved by Excel).
$ws.Range("A19:C20").EntireRow.Delete()

# --- 2. Row 17 loses its special "highlight" formatting ---------------------
# Before the edit, row 17 held "Break" and used the same highlighted style as
# the Lunch row. After the edit it is an ordinary schedule row, so copy the
# plain formatting from row 16 (a normal, un-highlighted row) onto row 17.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Rewrite the TIME / TOPIC text for the schedule table ----------------
$ws.Range("A6").Value = "9:00 - 9:30"
$ws.Range("C6").Value = "Apache Spark History/Background and MapReduce"

$ws.Range("A7").Value = "9:30 - 9:45"
$ws.Range("C7").Value = "System Setup"

$ws.Range("A8").Value = "9:45 - 10:15"
$ws.Range("C8").Value = "Apache Spark Context and RDDs"

$ws.Range("A9").Value = "10:15 - 10:45"
$ws.Range("C9").Value = "Exercise 1 - Running Spark Jobs"

$ws.Range("A10").Value = "10:45 - 11:15"
$ws.Range("C10").Value = "Apache Spark APIs"

$ws.Range("A11").Value = "11:15 - 11:45"
$ws.Range("C11").Value = "Exercise 2 - Access Logs"

$ws.Range("A12").Value = "11:45 - 1:00"
$ws.Range("C12").Value = "Lunch"

$ws.Range("A13").Value = "1:00 - 1:30"
$ws.Range("C13").Value = "Advanced Apache Spark APIs and Lineage"

$ws.Range("A14").Value = "1:30 - 2:00"
$ws.Range("C14").Value = "Exercise 3 - Joining Datasets"

$ws.Range("A15").Value = "2:00 - 2:30"
$ws.Range("C15").Value = "Shared Variables"

$ws.Range("A16").Value = "2:30 - 3:00"
$ws.Range("C16").Value = "Exercise 4 - Shared Variables"

$ws.Range("A17").Value = "3:00 - 3:15"
$ws.Range("C17").Value = "Misc. Concepts"

$ws.Range("A18").Value = "4:00 - 4:30"
$ws.Range("C18").Value = "Q&A"

# --- 4. Restore selection to C17, matching the saved workbook state --------
$ws.Range("C17").Select()
